$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44181
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 4000

# Row 3
$ws.Range("D3").Value = 44186

# Row 5
$ws.Range("D5").Value = 44179
$ws.Range("M5").Value = 45

# Row 6
$ws.Range("D6").Value = 44188
$ws.Range("M6").Value = 30

# Row 7
$ws.Range("D7").Value = 44193
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 3000
